# Update the appendix title:
#   "Experience shapes individual foraging specialization and success
#    in a virtual predator-prey system: Appendix 1"
# becomes
#   "Individual foraging specialization and success change across
#    experience in a virtual predator-prey system: Appendix 1"

$d = $word.ActiveDocument

# Drop the leading "Experience shapes " and capitalize "individual" -> "Individual"
$d.Content.Find.Execute(
    "Experience shapes individual", $true, $false, $false, $false, $false,
    $true, 1, $false, "Individual", 2
)

# Insert "change across experience" right after "success" (before "in a virtual")
$d.Content.Find.Execute(
    "success in a virtual", $true, $false, $false, $false, $false,
    $true, 1, $false, "success change across experience in a virtual", 2
)
